$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.456.01"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.568.65"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.25"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.23"
$ws.Range("E8").Value = "  -3.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3323"
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.134"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07475"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.68"
$ws.Range("E13").Value = "  -2.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.942"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.916"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").Value = "1.568.23"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001114"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06754"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.78"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.350"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.43"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("D24").Value = "22.450.31"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.386"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.567"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.94"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.70"
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.021"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.07"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").Value = "1.738.06"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.056"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.010"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.124"
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.684"
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08297"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02458"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2276"
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06397"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.371"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.288"
$ws.Range("E41").Value = "  -4.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.28"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6287"
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.82"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6130"
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.771"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.046"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.28"
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07231"
$ws.Range("E51").Value = "  -0.98%  "
